# Update "想去人数" (column F) figures on the "展览" and "全部类型" sheets.
# Same set of events is listed on both sheets; row numbers line up except
# that "全部类型" has a handful of extra rows interleaved, so each change
# below carries the row index for each sheet independently.

$wb = $excel.ActiveWorkbook

$sheet1Name = "展览"
$sheet4Name = "全部类型"

$ws1 = $wb.Worksheets.Item($sheet1Name)
$ws4 = $wb.Worksheets.Item($sheet4Name)

# Each entry: row on "展览", row on "全部类型", expected old value, new value.
$changes = @(
    @{ S1Row = 3;  S4Row = 3;  Old = 337;  New = 339 },
    @{ S1Row = 4;  S4Row = 4;  Old = 424;  New = 427 },
    @{ S1Row = 6;  S4Row = 6;  Old = 84;   New = 85 },
    @{ S1Row = 7;  S4Row = 7;  Old = 2176; New = 2177 },
    @{ S1Row = 10; S4Row = 10; Old = 106;  New = 107 },
    @{ S1Row = 11; S4Row = 11; Old = 4893; New = 4905 },
    @{ S1Row = 17; S4Row = 17; Old = 178;  New = 180 },
    @{ S1Row = 21; S4Row = 21; Old = 3839; New = 3853 },
    @{ S1Row = 22; S4Row = 22; Old = 706;  New = 709 },
    @{ S1Row = 23; S4Row = 23; Old = 648;  New = 657 },
    @{ S1Row = 26; S4Row = 26; Old = 102;  New = 103 },
    @{ S1Row = 27; S4Row = 27; Old = 116;  New = 117 },
    @{ S1Row = 28; S4Row = 28; Old = 21;   New = 22 },
    @{ S1Row = 30; S4Row = 30; Old = 86;   New = 88 },
    @{ S1Row = 32; S4Row = 32; Old = 8;    New = 9 },
    @{ S1Row = 34; S4Row = 35; Old = 929;  New = 939 },
    @{ S1Row = 35; S4Row = 36; Old = 2443; New = 2448 }
)

$mismatches = 0
foreach ($chg in $changes) {
    $cell1 = $ws1.Cells.Item($chg.S1Row, 6)
    if ($cell1.Value2 -ne $chg.Old) { $mismatches++ }
    $cell1.Value = $chg.New

    $cell4 = $ws4.Cells.Item($chg.S4Row, 6)
    if ($cell4.Value2 -ne $chg.Old) { $mismatches++ }
    $cell4.Value = $chg.New
}

Write-Output "Updated column F counts on sheets '$sheet1Name' and '$sheet4Name' ($($changes.Count) rows each, $mismatches unexpected prior value(s))."
